$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 562.6667
$ws.Range("I12").Value = 480.57144
$ws.Range("J12").Value = 677.6
$ws.Range("K12").Value = 480.57144
$ws.Range("L12").Value = 677.6
$ws.Range("M12").Value = -310.57144
$ws.Range("N12").Value = -1017.6
$ws.Range("H19").Value = 1536.7858
$ws.Range("J19").Value = 2002.4445
$ws.Range("L19").Value = 2002.4445
$ws.Range("N19").Value = -2352.4445
$ws.Range("H38").Value = 1454.5333
$ws.Range("I38").Value = 1005.1667
$ws.Range("J38").Value = 1754.1111
$ws.Range("K38").Value = 3015.5001
$ws.Range("L38").Value = 5262.3333
$ws.Range("M38").Value = -2643.5001
$ws.Range("N38").Value = -6006.3333
$ws.Range("H41").Value = 416.25
$ws.Range("I41").Value = 350
$ws.Range("J41").Value = 456
$ws.Range("K41").Value = 350
$ws.Range("L41").Value = 456
$ws.Range("M41").Value = 90
$ws.Range("N41").Value = -1336
$ws.Range("H86").Value = 5250
$ws.Range("I86").Value = 3500
$ws.Range("J86").Value = 7000
$ws.Range("K86").Value = 3500
$ws.Range("L86").Value = 7000
$ws.Range("M86").Value = -2377
$ws.Range("N86").Value = -9246
$ws.Range("H89").Value = 5250
$ws.Range("I89").Value = 3500
$ws.Range("J89").Value = 7000
$ws.Range("K89").Value = 17500
$ws.Range("L89").Value = 35000
$ws.Range("M89").Value = -11884
$ws.Range("N89").Value = -46232
$ws.Range("H116").Value = 9978.929
$ws.Range("J116").Value = 3178.875
$ws.Range("L116").Value = 3178.875
$ws.Range("N116").Value = -10062.875
$ws.Range("H121").Value = 1382.7693
$ws.Range("J121").Value = 1382.7693
$ws.Range("L121").Value = 4148.3079
$ws.Range("N121").Value = -7642.3079
$ws.Range("H124").Value = 46102.95
$ws.Range("J124").Value = 46102.95
$ws.Range("L124").Value = 46102.95
$ws.Range("N124").Value = -55922.95
$ws.Range("H137").Value = 2118.6553
$ws.Range("I137").Value = 1974.5
$ws.Range("K137").Value = 5923.5
$ws.Range("M137").Value = -3373.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1641.1428
$ws.Range("I2").Value = 1331.3334
$ws.Range("K2").Value = 1331.3334
$ws.Range("M2").Value = -1218.3334
$ws.Range("H32").Value = 3407
$ws.Range("I32").Value = 2474.8865
$ws.Range("J32").Value = 6336.5
$ws.Range("K32").Value = 2474.8865
$ws.Range("L32").Value = 6336.5
$ws.Range("M32").Value = -2187.8865
$ws.Range("N32").Value = -6910.5
$ws.Range("H61").Value = 2563.9285
$ws.Range("I61").Value = 1811.5385
$ws.Range("K61").Value = 1811.5385
$ws.Range("M61").Value = -1599.5385
$ws.Range("H74").Value = 1680.1111
$ws.Range("I74").Value = 1579.125
$ws.Range("K74").Value = 1579.125
$ws.Range("M74").Value = -705.125
$ws.Range("H77").Value = 1680.1111
$ws.Range("I77").Value = 1579.125
$ws.Range("K77").Value = 7895.625
$ws.Range("M77").Value = -3527.625
$ws.Range("H116").Value = 1641.1428
$ws.Range("I116").Value = 1331.3334
$ws.Range("K116").Value = 1331.3334
$ws.Range("M116").Value = 962.6666
$ws.Range("H122").Value = 1725.8572
$ws.Range("I122").Value = 1752.15
$ws.Range("K122").Value = 5256.450000000001
$ws.Range("M122").Value = -2806.450000000001
$ws.Range("H132").Value = 4149
$ws.Range("I132").Value = 3965.6667
$ws.Range("J132").Value = 4332.3335
$ws.Range("K132").Value = 11897.0001
$ws.Range("L132").Value = 12997.0005
$ws.Range("M132").Value = -9367.000100000001
$ws.Range("N132").Value = -18057.0005
$ws.Range("H136").Value = 2563.9285
$ws.Range("I136").Value = 1811.5385
$ws.Range("K136").Value = 5434.6155
$ws.Range("M136").Value = -2884.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1641.1428
$ws.Range("I3").Value = 1331.3334
$ws.Range("K3").Value = 1331.3334
$ws.Range("M3").Value = -1217.3334
$ws.Range("H94").Value = 672.73334
$ws.Range("I94").Value = 448.3
$ws.Range("J94").Value = 1121.6
$ws.Range("K94").Value = 448.3
$ws.Range("L94").Value = 1121.6
$ws.Range("M94").Value = 2.699999999999989
$ws.Range("N94").Value = -2023.6
$ws.Range("H95").Value = 71893.25
$ws.Range("J95").Value = 71893.25
$ws.Range("L95").Value = 71893.25
$ws.Range("N95").Value = -77385.25
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
$ws.Range("H99").Value = 1234.625
$ws.Range("I99").Value = 1121.3334
$ws.Range("K99").Value = 1121.3334
$ws.Range("M99").Value = 376.6666
$ws.Range("H134").Value = 4363.4814
$ws.Range("I134").Value = 4662.0527
$ws.Range("K134").Value = 13986.1581
$ws.Range("M134").Value = -11451.1581

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 160
$ws.Range("I7").Value = 75
$ws.Range("K7").Value = 75
$ws.Range("M7").Value = 38
$ws.Range("H16").Value = 835.44446
$ws.Range("I16").Value = 845.1429000000001
$ws.Range("J16").Value = 801.5
$ws.Range("K16").Value = 845.1429000000001
$ws.Range("L16").Value = 801.5
$ws.Range("M16").Value = -558.1429000000001
$ws.Range("N16").Value = -1375.5
$ws.Range("H28").Value = 40643
$ws.Range("J28").Value = 40643
$ws.Range("L28").Value = 40643
$ws.Range("N28").Value = -41133
$ws.Range("H31").Value = 4701.077
$ws.Range("I31").Value = 1806
$ws.Range("J31").Value = 5987.778
$ws.Range("K31").Value = 1806
$ws.Range("L31").Value = 5987.778
$ws.Range("M31").Value = -1511
$ws.Range("N31").Value = -6577.778
$ws.Range("H34").Value = 4701.077
$ws.Range("I34").Value = 1806
$ws.Range("J34").Value = 5987.778
$ws.Range("K34").Value = 1806
$ws.Range("L34").Value = 5987.778
$ws.Range("M34").Value = -1604
$ws.Range("N34").Value = -6391.778
$ws.Range("H107").Value = 369.29413
$ws.Range("I107").Value = 317.14285
$ws.Range("J107").Value = 612.6667
$ws.Range("K107").Value = 317.14285
$ws.Range("L107").Value = 612.6667
$ws.Range("M107").Value = 1602.85715
$ws.Range("N107").Value = -4452.6667
$ws.Range("H113").Value = 835.44446
$ws.Range("I113").Value = 845.1429000000001
$ws.Range("J113").Value = 801.5
$ws.Range("K113").Value = 845.1429000000001
$ws.Range("L113").Value = 801.5
$ws.Range("M113").Value = 1324.8571
$ws.Range("N113").Value = -5141.5
$ws.Range("H122").Value = 1867.8572
$ws.Range("I122").Value = 1959.1111
$ws.Range("K122").Value = 5877.3333
$ws.Range("M122").Value = -3427.3333
$ws.Range("H134").Value = 3364.1
$ws.Range("I134").Value = 2580.25
$ws.Range("K134").Value = 7740.75
$ws.Range("M134").Value = -5205.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3688.6667
$ws.Range("I3").Value = 1900
$ws.Range("J3").Value = 4583
$ws.Range("K3").Value = 5700
$ws.Range("L3").Value = 13749
$ws.Range("M3").Value = -5588
$ws.Range("N3").Value = -13973
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H122").Value = 1278.909
$ws.Range("I122").Value = 1016.8
$ws.Range("K122").Value = 9151.199999999999
$ws.Range("M122").Value = -6701.199999999999
$ws.Range("H137").Value = 3247.9656
$ws.Range("J137").Value = 4736.706
$ws.Range("L137").Value = 14210.118
$ws.Range("N137").Value = -24410.118

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 20000000
$ws.Range("J47").Value = 20000000
$ws.Range("L47").Value = 20000000
$ws.Range("N47").Value = -20001136
$ws.Range("H48").Value = 1027
$ws.Range("I48").Value = 1027
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 1027
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -542
$ws.Range("N48").ClearContents()
$ws.Range("H102").Value = 2122.4736
$ws.Range("I102").Value = 2611.2222
$ws.Range("K102").Value = 2611.2222
$ws.Range("M102").Value = -989.2222000000002
$ws.Range("H136").Value = 6749.5
$ws.Range("J136").Value = 6749.5
$ws.Range("L136").Value = 20248.5
$ws.Range("N136").Value = -25348.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1540.8182
$ws.Range("I46").Value = 600
$ws.Range("K46").Value = 600
$ws.Range("M46").Value = -412
$ws.Range("H48").Value = 8041
$ws.Range("I48").Value = 8041
$ws.Range("K48").Value = 8041
$ws.Range("M48").Value = -7380
$ws.Range("H132").Value = 3339.4
$ws.Range("I132").Value = 2749
$ws.Range("K132").Value = 8247
$ws.Range("M132").Value = -5717

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 70049
$ws.Range("J42").Value = 70049
$ws.Range("L42").Value = 70049
$ws.Range("N42").Value = -70805
$ws.Range("H51").Value = 38888
$ws.Range("J51").Value = 38888
$ws.Range("L51").Value = 38888
$ws.Range("N51").Value = -39908
$ws.Range("H52").Value = 18749.5
$ws.Range("I52").Value = 9999
$ws.Range("K52").Value = 9999
$ws.Range("M52").Value = -9773
$ws.Range("H54").Value = 28933.334
$ws.Range("J54").Value = 28933.334
$ws.Range("L54").Value = 28933.334
$ws.Range("N54").Value = -29973.334
$ws.Range("H122").Value = 47130.883
$ws.Range("I122").Value = 56913.285
$ws.Range("J122").Value = 1479.6666
$ws.Range("K122").Value = 170739.855
$ws.Range("L122").Value = 4438.9998
$ws.Range("M122").Value = -168289.855
$ws.Range("N122").Value = -9338.9998
$ws.Range("H132").Value = 2657.8333
$ws.Range("I132").Value = 1979.6
$ws.Range("J132").Value = 3142.2856
$ws.Range("K132").Value = 5938.799999999999
$ws.Range("L132").Value = 9426.856800000001
$ws.Range("M132").Value = -3408.799999999999
$ws.Range("N132").Value = -14486.8568
